$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip the sign of J14:J25 (1 -> -1)
$ws.Range("J14:J25").Value = -1

# Update the active selection to J28 (matches final sheetView selection in diff)
$ws.Range("J28").Select()
